$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "27.857.56"
$ws.Cells.Item(2, 5).Value = "  +0.56%  "

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "1.766.43"
$ws.Cells.Item(3, 5).Value = "  +0.51%  "

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) "1.002"
$ws.Cells.Item(4, 5).Value = "  +0.04%  "

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "327.32"
$ws.Cells.Item(5, 5).Value = "  +0.76%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +0.08%  "

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "0.4476"
$ws.Cells.Item(7, 5).Value = "  -3.06%  "

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "0.3540"
$ws.Cells.Item(8, 5).Value = "  -1.84%  "

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "0.07436"
$ws.Cells.Item(9, 5).Value = "  -1.07%  "

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "42.06"
$ws.Cells.Item(10, 5).Value = "  -0.34%  "

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "1.093"
$ws.Cells.Item(11, 5).Value = "  -0.56%  "

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "1.001"
$ws.Cells.Item(12, 5).Value = "  +0.03%  "

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "20.81"
$ws.Cells.Item(13, 5).Value = "  +0.09%  "

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "6.017"

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "7.194"
$ws.Cells.Item(15, 5).Value = "  +0.99%  "

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "1.769.66"
$ws.Cells.Item(16, 5).Value = "  +0.69%  "

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "92.87"
$ws.Cells.Item(17, 5).Value = "  +0.57%  "

# Row 18
$ws.Cells.Item(18, 5).Value = "  -1.03%  "

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) "0.06427"
$ws.Cells.Item(19, 5).Value = "  +0.44%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +0.09%  "

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "17.13"
$ws.Cells.Item(21, 5).Value = "  +1.94%  "

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "5.782"
$ws.Cells.Item(22, 5).Value = "  -0.50%  "

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) "27.922.04"
$ws.Cells.Item(23, 5).Value = "  +0.57%  "

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "11.26"
$ws.Cells.Item(24, 5).Value = "  -0.14%  "

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "2.125"
$ws.Cells.Item(25, 5).Value = "  +0.62%  "

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "162.63"
$ws.Cells.Item(26, 5).Value = "  -1.05%  "

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) "20.16"
$ws.Cells.Item(27, 5).Value = "  -1.10%  "

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) "1.975.92"
$ws.Cells.Item(28, 5).Value = "  +0.63%  "

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) "2.173"
$ws.Cells.Item(29, 5).Value = "  +4.46%  "

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) "124.88"
$ws.Cells.Item(30, 5).Value = "  -1.23%  "

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) "1.098"
$ws.Cells.Item(31, 5).Value = "  +3.43%  "

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) "0.09141"
$ws.Cells.Item(32, 5).Value = "  -1.30%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "HuobiToken"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Cells.Item(33, 4) "3.648"
$ws.Cells.Item(33, 5).Value = "  -0.60%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "Filecoin"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Cells.Item(34, 4) "5.562"
$ws.Cells.Item(34, 5).Value = "  +0.52%  "

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) "11.83"
$ws.Cells.Item(35, 5).Value = "  -0.65%  "

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) "0.02291"
$ws.Cells.Item(36, 5).Value = "  -0.49%  "

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) "0.06098"
$ws.Cells.Item(37, 5).Value = "  +0.98%  "

# Row 38
Set-TextValue $ws.Cells.Item(38, 4) "0.2090"
$ws.Cells.Item(38, 5).Value = "  -0.59%  "

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "4.970"
$ws.Cells.Item(39, 5).Value = "  -0.10%  "

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "0.6299"
$ws.Cells.Item(40, 5).Value = "  -0.87%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -1.75%  "

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "1.387"
$ws.Cells.Item(42, 5).Value = "  +0.53%  "

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "7.928"
$ws.Cells.Item(43, 5).Value = "  +1.18%  "

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "13.31"
$ws.Cells.Item(44, 5).Value = "  -0.06%  "

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "3.734"
$ws.Cells.Item(45, 5).Value = "  +0.47%  "

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "0.5843"
$ws.Cells.Item(46, 5).Value = "  -1.23%  "

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "122.01"
$ws.Cells.Item(47, 5).Value = "  -1.16%  "

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) "1.947"
$ws.Cells.Item(48, 5).Value = "  -0.28%  "

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "0.06902"

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) "1.135"
$ws.Cells.Item(50, 5).Value = "  -1.28%  "

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "72.72"
$ws.Cells.Item(51, 5).Value = "  +0.49%  "
